$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "15÷3=5, 0"
$t.Cell(1, 2).Range.Text = "77÷2=38, 1"
$t.Cell(1, 3).Range.Text = "43÷8=5, 3"
$t.Cell(1, 4).Range.Text = "98÷8=12, 2"
$t.Cell(1, 5).Range.Text = "42÷8=5, 2"
$t.Cell(5, 1).Range.Text = "69÷6=11, 3"
$t.Cell(5, 2).Range.Text = "12÷3=4, 0"
$t.Cell(5, 3).Range.Text = "92÷4=23, 0"
$t.Cell(5, 4).Range.Text = "58÷6=9, 4"
$t.Cell(5, 5).Range.Text = "86÷4=21, 2"
$t.Cell(9, 1).Range.Text = "59÷4=14, 3"
$t.Cell(9, 2).Range.Text = "85÷5=17, 0"
$t.Cell(9, 3).Range.Text = "47÷2=23, 1"
$t.Cell(9, 4).Range.Text = "38÷7=5, 3"
$t.Cell(9, 5).Range.Text = "18÷7=2, 4"
$t.Cell(13, 1).Range.Text = "15÷6=2, 3"
$t.Cell(13, 2).Range.Text = "36÷2=18, 0"
$t.Cell(13, 3).Range.Text = "59÷4=14, 3"
$t.Cell(13, 4).Range.Text = "42÷7=6, 0"
$t.Cell(13, 5).Range.Text = "21÷9=2, 3"
$t.Cell(17, 1).Range.Text = "11÷8=1, 3"
$t.Cell(17, 2).Range.Text = "36÷7=5, 1"
$t.Cell(17, 3).Range.Text = "36÷3=12, 0"
$t.Cell(17, 4).Range.Text = "18÷7=2, 4"
$t.Cell(17, 5).Range.Text = "92÷9=10, 2"
